$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formatting Mobile Numbers")

# New static values to replace the RANDBETWEEN shared-formula results
$values = @{
    2  = 4759462429
    3  = 3111843265
    4  = 3007723039
    5  = 1530049582
    6  = 9584299252
    7  = 7720819933
    8  = 3447612857
    9  = 9486032302
    10 = 8062846356
    11 = 2279190723
    12 = 8371665677
    13 = 5265926213
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

$ws.Range("E5").Select()
